$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply an AutoFilter on column L ("Employment Type", field 12, colId=11)
# keeping only rows where the value is "Consultant". Operator 7 = xlFilterValues
# (the "standard" / value-list filter that serialises as <filters><filter .../></filters>
# rather than the <customFilters> form).
$range = $ws.Range("A1:P109")
$range.AutoFilter(12, @("Consultant"), 7)

# Re-establish the frozen panes (unchanged: freeze header row + first two columns)
# and move the active selection to L87, matching the post-edit view state.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("C2").Select()
$win.FreezePanes = $true
$ws.Range("L87").Select()
